# Expense Tracker Application Report - apply commit "Ppt content almost done. brief report started"
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix "it's location" -> "its location" and append the new sentence about
#    file locations to the "brief report" intro paragraph.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("it’s location", $true, $false, $false, $false, $false, `
    $true, 1, $false, "its location", 2)

$introPara = $d.Paragraphs(2)
$introEnd = $introPara.Range
$introEnd.Collapse(0)
$introEnd.InsertAfter(" All files are within the /flaskDemo/ folder and further sub directories will be noted")

# ---------------------------------------------------------------------------
# 2) Split out "MySQLAlchemy" (cosmetic run split in the source diff; content
#    stays the same) - just make sure the sentence text is intact.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(" should use Flask with Flask-MySQLAlchemy.", $true, $false, $false, $false, $false, `
    $true, 1, $false, " should use Flask with Flask-MySQLAlchemy.", 2)

# ---------------------------------------------------------------------------
# 3) Populate the previously-empty cells in the "Python/MySQL/Flask" row with
#    the file locations, line info and achieved-task description.
# ---------------------------------------------------------------------------
$row2 = $d.Tables(1).Rows(2)

$fileCell = $row2.Cells(2)
$fileCell.Range.Text = "forms.py`rmodels.py`rroutes.py"

$lineCell = $row2.Cells(3)
$lineCell.Range.Text = "Entire File"

$achievedCell = $row2.Cells(4)
$achievedCell.Range.Text = "Allows connection between database, web application functionality, and deployment"

# ---------------------------------------------------------------------------
# 4) Remove the trailing "using" from the structural-change task description.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Include at least one structural change to the database using ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Include at least one structural change to the database ", 2)

# ---------------------------------------------------------------------------
# 5) Move "compound condition" from the end of the first run to the start of
#    the second run in the SELECT task description.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("Include at least one SELECT using a compound condition ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Include at least one SELECT using a ", 2)
$null = $d.Content.Find.Execute("using regular SQL, and also the equivalent of a compound condition select using Flask-SQLAlchemy.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "compound condition using regular SQL, and also the equivalent of a compound condition select using Flask-SQLAlchemy.", 2)
